# Apply edits described by the commit "Added many more features"
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# 1. Title heading + duplicated bold SEO title (both occurrences updated identically)
Replace-Text "Play Black Horse Deluxe for Free - Features Customizable Volatility" `
             "Play Black Horse Deluxe for Free - Exciting Wild West Slot Game"

# 2. "What we like" bullet list - rewording / shuffled content
Replace-Text "Customizable volatility levels." "Customizable volatility level"
Replace-Text "High payouts up to 1,300x the initial bet." "Free Spins feature and Gamble option"
Replace-Text "Free Spins feature and Gamble option." "High payouts of up to 1,300 times the bet"
Replace-Text "Unique Wild West theme and symbols." "Unique Wild West theme"

# 3. "What we don't like" bullet list
Replace-Text "Not the most original slot." "Not the most original title"
Replace-Text "May not appeal to players who are not interested in Western themes." "Limited number of paylines"

# 4. Meta description (italic paragraph)
Replace-Text "Read our review of Black Horse Deluxe, a unique online slot game with Wild West symbols. Play for free and enjoy customizable volatility and high payouts." `
             "Read our review of Black Horse Deluxe, a Wild West themed slot game with high payouts. Play for free and enjoy the exciting features."
